$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Clec11a"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.059611
$ws.Cells.Item(2, 8).Value = 0.178833
$ws.Cells.Item(2, 9).Value = 0.003943014985542741
$ws.Cells.Item(2, 10).Value = 0.003943014985542741
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 168.1098273333333
$ws.Cells.Item(2, 14).Value = 504.329482
$ws.Cells.Item(2, 15).Value = 0.2984182258032519
$ws.Cells.Item(2, 16).Value = 0.298418225803252
$ws.Cells.Item(2, 17).Value = 10.02119491716733
$ws.Cells.Item(2, 18).Value = 90.19075425450599
$ws.Cells.Item(2, 19).Value = 0.0011766675363013
$ws.Cells.Item(2, 20).Value = 0.0011766675363013

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Clec11a"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.059611
$ws.Cells.Item(3, 8).Value = 0.178833
$ws.Cells.Item(3, 9).Value = 0.003943014985542741
$ws.Cells.Item(3, 10).Value = 0.003943014985542741
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 163.0062356666667
$ws.Cells.Item(3, 14).Value = 489.018707
$ws.Cells.Item(3, 15).Value = 0.2893586437755394
$ws.Cells.Item(3, 16).Value = 0.2893586437755394
$ws.Cells.Item(3, 17).Value = 9.716964714325666
$ws.Cells.Item(3, 18).Value = 87.452682428931
$ws.Cells.Item(3, 19).Value = 0.001140945468603276
$ws.Cells.Item(3, 20).Value = 0.001140945468603276

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Clec11a"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.059611
$ws.Cells.Item(4, 8).Value = 0.178833
$ws.Cells.Item(4, 9).Value = 0.003943014985542741
$ws.Cells.Item(4, 10).Value = 0.003943014985542741
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 165.99353
$ws.Cells.Item(4, 14).Value = 497.98059
$ws.Cells.Item(4, 15).Value = 0.294661504941043
$ws.Cells.Item(4, 16).Value = 0.294661504941043
$ws.Cells.Item(4, 17).Value = 9.895040316829999
$ws.Cells.Item(4, 18).Value = 89.05536285146999
$ws.Cells.Item(4, 19).Value = 0.001161854729645109
$ws.Cells.Item(4, 20).Value = 0.001161854729645109

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Clec11a"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.059611
$ws.Cells.Item(5, 8).Value = 0.178833
$ws.Cells.Item(5, 9).Value = 0.003943014985542741
$ws.Cells.Item(5, 10).Value = 0.003943014985542741
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 66.22673433333334
$ws.Cells.Item(5, 14).Value = 198.680203
$ws.Cells.Item(5, 15).Value = 0.1175616254801657
$ws.Cells.Item(5, 16).Value = 0.1175616254801657
$ws.Cells.Item(5, 17).Value = 3.947841860344333
$ws.Cells.Item(5, 18).Value = 35.530576743099
$ws.Cells.Item(5, 19).Value = 0.0004635472509930566
$ws.Cells.Item(5, 20).Value = 0.0004635472509930566

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Clec11a"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 14.75172666666667
$ws.Cells.Item(6, 8).Value = 44.25518
$ws.Cells.Item(6, 9).Value = 0.9757641930062764
$ws.Cells.Item(6, 10).Value = 0.9757641930062765
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 168.1098273333333
$ws.Cells.Item(6, 14).Value = 504.329482
$ws.Cells.Item(6, 15).Value = 0.2984182258032519
$ws.Cells.Item(6, 16).Value = 0.298418225803252
$ws.Cells.Item(6, 17).Value = 2479.910222801862
$ws.Cells.Item(6, 18).Value = 22319.19200521676
$ws.Cells.Item(6, 19).Value = 0.2911858192792749
$ws.Cells.Item(6, 20).Value = 0.291185819279275

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Clec11a"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 14.75172666666667
$ws.Cells.Item(7, 8).Value = 44.25518
$ws.Cells.Item(7, 9).Value = 0.9757641930062764
$ws.Cells.Item(7, 10).Value = 0.9757641930062765
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 163.0062356666667
$ws.Cells.Item(7, 14).Value = 489.018707
$ws.Cells.Item(7, 15).Value = 0.2893586437755394
$ws.Cells.Item(7, 16).Value = 0.2893586437755394
$ws.Cells.Item(7, 17).Value = 2404.623433516918
$ws.Cells.Item(7, 18).Value = 21641.61090165226
$ws.Cells.Item(7, 19).Value = 0.2823458035330298
$ws.Cells.Item(7, 20).Value = 0.2823458035330298

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Clec11a"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.75172666666667
$ws.Cells.Item(8, 8).Value = 44.25518
$ws.Cells.Item(8, 9).Value = 0.9757641930062764
$ws.Cells.Item(8, 10).Value = 0.9757641930062765
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 165.99353
$ws.Cells.Item(8, 14).Value = 497.98059
$ws.Cells.Item(8, 15).Value = 0.294661504941043
$ws.Cells.Item(8, 16).Value = 0.294661504941043
$ws.Cells.Item(8, 17).Value = 2448.691182995134
$ws.Cells.Item(8, 18).Value = 22038.2206469562
$ws.Cells.Item(8, 19).Value = 0.2875201455788117
$ws.Cells.Item(8, 20).Value = 0.2875201455788118

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Clec11a"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.75172666666667
$ws.Cells.Item(9, 8).Value = 44.25518
$ws.Cells.Item(9, 9).Value = 0.9757641930062764
$ws.Cells.Item(9, 10).Value = 0.9757641930062765
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 66.22673433333334
$ws.Cells.Item(9, 14).Value = 198.680203
$ws.Cells.Item(9, 15).Value = 0.1175616254801657
$ws.Cells.Item(9, 16).Value = 0.1175616254801657
$ws.Cells.Item(9, 17).Value = 976.9586829112824
$ws.Cells.Item(9, 18).Value = 8792.628146201541
$ws.Cells.Item(9, 19).Value = 0.1147124246151599
$ws.Cells.Item(9, 20).Value = 0.11471242461516

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Clec11a"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.306789
$ws.Cells.Item(10, 8).Value = 0.9203669999999999
$ws.Cells.Item(10, 9).Value = 0.0202927920081809
$ws.Cells.Item(10, 10).Value = 0.02029279200818091
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 168.1098273333333
$ws.Cells.Item(10, 14).Value = 504.329482
$ws.Cells.Item(10, 15).Value = 0.2984182258032519
$ws.Cells.Item(10, 16).Value = 0.298418225803252
$ws.Cells.Item(10, 17).Value = 51.57424581776599
$ws.Cells.Item(10, 18).Value = 464.168212359894
$ws.Cells.Item(10, 19).Value = 0.006055738987675756
$ws.Cells.Item(10, 20).Value = 0.006055738987675758

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Clec11a"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.306789
$ws.Cells.Item(11, 8).Value = 0.9203669999999999
$ws.Cells.Item(11, 9).Value = 0.0202927920081809
$ws.Cells.Item(11, 10).Value = 0.02029279200818091
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 163.0062356666667
$ws.Cells.Item(11, 14).Value = 489.018707
$ws.Cells.Item(11, 15).Value = 0.2893586437755394
$ws.Cells.Item(11, 16).Value = 0.2893586437755394
$ws.Cells.Item(11, 17).Value = 50.008520033941
$ws.Cells.Item(11, 18).Value = 450.076680305469
$ws.Cells.Item(11, 19).Value = 0.005871894773906332
$ws.Cells.Item(11, 20).Value = 0.005871894773906333

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Clec11a"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.306789
$ws.Cells.Item(12, 8).Value = 0.9203669999999999
$ws.Cells.Item(12, 9).Value = 0.0202927920081809
$ws.Cells.Item(12, 10).Value = 0.02029279200818091
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 165.99353
$ws.Cells.Item(12, 14).Value = 497.98059
$ws.Cells.Item(12, 15).Value = 0.294661504941043
$ws.Cells.Item(12, 16).Value = 0.294661504941043
$ws.Cells.Item(12, 17).Value = 50.92498907517
$ws.Cells.Item(12, 18).Value = 458.32490167653
$ws.Cells.Item(12, 19).Value = 0.005979504632586155
$ws.Cells.Item(12, 20).Value = 0.005979504632586156

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Clec11a"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.306789
$ws.Cells.Item(13, 8).Value = 0.9203669999999999
$ws.Cells.Item(13, 9).Value = 0.0202927920081809
$ws.Cells.Item(13, 10).Value = 0.02029279200818091
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 66.22673433333334
$ws.Cells.Item(13, 14).Value = 198.680203
$ws.Cells.Item(13, 15).Value = 0.1175616254801657
$ws.Cells.Item(13, 16).Value = 0.1175616254801657
$ws.Cells.Item(13, 17).Value = 20.317633599389
$ws.Cells.Item(13, 18).Value = 182.858702394501
$ws.Cells.Item(13, 19).Value = 0.002385653614012662
$ws.Cells.Item(13, 20).Value = 0.002385653614012662

